$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new export timestamp (2024-12-03 -> 2024-12-04)
$ws.Name = "IClientBalance-20241204-100653-"

# Bump the date column (G) for every data row (2-274) from 2024-12-03 (45629)
# to 2024-12-04 (45630), keeping the existing date number formatting.
$ws.Range("G2:G274").Value = 45630
